$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.328.78"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.14"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.37"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4712"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2901"
$ws.Range("E8").Value = "  +2.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06622"

$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08020"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.35"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.878.42"
$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.159"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6900"
$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.39"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.314.43"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.04"
$ws.Range("E18").Value = "  +5.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007733"
$ws.Range("E19").Value = "  +5.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.122.58"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.318"
$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.76"
$ws.Range("E25").Value = "  +0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.289"
$ws.Range("E26").Value = "  +1.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.00"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.965"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09952"
$ws.Range("E30").Value = "  +2.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.372"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.464"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.096"
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04706"
$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7036"
$ws.Range("E36").Value = "  -0.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01883"
$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.643"
$ws.Range("E39").Value = "  +2.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.316"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.38"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4174"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.63"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.149"
$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.233"
$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "933.72"
$ws.Range("E49").Value = "  -3.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.51"
$ws.Range("E50").Value = "  +1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05667"
$ws.Range("E51").Value = "  +0.39%  "
